# error handle in login
# Adds header label cells E1:I1 containing the text "0","1","2","3","4"
# (stored as genuine text, not numbers) and extends the used range down
# to an empty row 6, matching the target worksheet shape (A1:I6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write "0".."4" into E1:I1 as real text values (not numbers).
# We build each value via a temporary text-formula cell and paste only
# the resulting value back in, which yields a plain text cell without
# leaving a formula behind and without forcing a "quote prefix" style.
$helper = $ws.Cells.Item(1, 20)
for ($i = 0; $i -lt 5; $i++) {
    $col = 5 + $i
    $helper.Formula = "=""" + $i + """"
    $helper.Copy()
    $ws.Cells.Item(1, $col).PasteSpecial(-4163)  # xlPasteValues
}
$helper.ClearContents()
$excel.CutCopyMode = $false

# Touch row 6 so it becomes part of the sheet's used range (empty row),
# pushing the worksheet dimension down to A1:I6.
$ws.Cells.Item(6, 1).Style = "Normal"
